$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 146.66667  # H9: 168.6 -> 146.66667
$ws.Cells.Item(9, 9).Value = 100  # I9: 250.5 -> 100
$ws.Cells.Item(9, 10).Value = 170  # J9: 114 -> 170
$ws.Cells.Item(9, 11).Value = 100  # K9: 250.5 -> 100
$ws.Cells.Item(9, 12).Value = 170  # L9: 114 -> 170
$ws.Cells.Item(9, 13).Value = 69  # M9: -81.5 -> 69
$ws.Cells.Item(9, 14).Value = -508  # N9: -452 -> -508
$ws.Cells.Item(19, 8).Value = 5699.6665  # H19: 6699.4 -> 5699.6665
$ws.Cells.Item(19, 9).Value = 6639.8  # I19: 8124.5 -> 6639.8
$ws.Cells.Item(19, 11).Value = 6639.8  # K19: 8124.5 -> 6639.8
$ws.Cells.Item(19, 13).Value = -6464.8  # M19: -7949.5 -> -6464.8
$ws.Cells.Item(32, 8).Value = 2000  # H32: 966.6667 -> 2000
$ws.Cells.Item(32, 10).Value = 2000  # J32: 966.6667 -> 2000
$ws.Cells.Item(32, 12).Value = 2000  # L32: 966.6667 -> 2000
$ws.Cells.Item(32, 14).Value = -2652  # N32: -1618.6667 -> -2652
$ws.Cells.Item(43, 8).Value = 999  # H43: 0 -> 999
$ws.Cells.Item(43, 9).Value = 999  # I43: 0 -> 999
$ws.Cells.Item(43, 11).Value = 999  # K43: 0 -> 999
$ws.Cells.Item(43, 13).Value = -930  # M43: None -> -930
$ws.Cells.Item(86, 8).Value = 2000  # H86: 1949.5 -> 2000
$ws.Cells.Item(86, 9).Value = 0  # I86: 1899 -> 0
$ws.Cells.Item(86, 11).Value = 0  # K86: 1899 -> 0
$ws.Cells.Item(86, 13).ClearContents()  # M86: remove (was -776)
$ws.Cells.Item(89, 8).Value = 2000  # H89: 1949.5 -> 2000
$ws.Cells.Item(89, 9).Value = 0  # I89: 1899 -> 0
$ws.Cells.Item(89, 11).Value = 0  # K89: 9495 -> 0
$ws.Cells.Item(89, 13).ClearContents()  # M89: remove (was -3879)
$ws.Cells.Item(107, 8).Value = 280.66666  # H107: 281.6 -> 280.66666
$ws.Cells.Item(107, 10).Value = 281  # J107: 283.5 -> 281
$ws.Cells.Item(107, 12).Value = 281  # L107: 283.5 -> 281
$ws.Cells.Item(107, 14).Value = -4121  # N107: -4123.5 -> -4121
$ws.Cells.Item(135, 8).Value = 226.875  # H135: 279.44446 -> 226.875
$ws.Cells.Item(135, 9).Value = 187.85715  # I135: 189.375 -> 187.85715
$ws.Cells.Item(135, 10).Value = 500  # J135: 1000 -> 500
$ws.Cells.Item(135, 11).Value = 1690.71435  # K135: 1704.375 -> 1690.71435
$ws.Cells.Item(135, 12).Value = 4500  # L135: 9000 -> 4500
$ws.Cells.Item(135, 13).Value = 844.28565  # M135: 830.625 -> 844.28565
$ws.Cells.Item(135, 14).Value = -9570  # N135: -14070 -> -9570
$ws.Cells.Item(137, 8).Value = 2400  # H137: 2666.6667 -> 2400
$ws.Cells.Item(138, 8).Value = 2641.9333  # H138: 2787.8572 -> 2641.9333
$ws.Cells.Item(138, 9).Value = 1359.5  # I138: 1468.1428 -> 1359.5
$ws.Cells.Item(138, 11).Value = 4078.5  # K138: 4404.428400000001 -> 4078.5
$ws.Cells.Item(138, 13).Value = 1061.5  # M138: 735.5715999999993 -> 1061.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2296.2046  # H32: 2420.805 -> 2296.2046
$ws.Cells.Item(32, 9).Value = 2036.5  # I32: 2147.513 -> 2036.5
$ws.Cells.Item(32, 11).Value = 2036.5  # K32: 2147.513 -> 2036.5
$ws.Cells.Item(32, 13).Value = -1749.5  # M32: -1860.513 -> -1749.5
$ws.Cells.Item(37, 8).Value = 13998  # H37: 9997.5 -> 13998
$ws.Cells.Item(37, 10).Value = 21999  # J37: 0 -> 21999
$ws.Cells.Item(37, 12).Value = 21999  # L37: 0 -> 21999
$ws.Cells.Item(37, 14).Value = -22545  # N37: None -> -22545
$ws.Cells.Item(45, 8).Value = 2453.8  # H45: 4000 -> 2453.8
$ws.Cells.Item(45, 9).Value = 2492.25  # I45: 4000 -> 2492.25
$ws.Cells.Item(45, 10).Value = 2300  # J45: 0 -> 2300
$ws.Cells.Item(45, 11).Value = 2492.25  # K45: 4000 -> 2492.25
$ws.Cells.Item(45, 12).Value = 2300  # L45: 0 -> 2300
$ws.Cells.Item(45, 13).Value = -2115.25  # M45: -3623 -> -2115.25
$ws.Cells.Item(45, 14).Value = -3054  # N45: None -> -3054
$ws.Cells.Item(55, 8).Value = 56053  # H55: 0 -> 56053
$ws.Cells.Item(55, 10).Value = 56053  # J55: 0 -> 56053
$ws.Cells.Item(55, 12).Value = 56053  # L55: 0 -> 56053
$ws.Cells.Item(55, 14).Value = -56683  # N55: None -> -56683
$ws.Cells.Item(61, 8).Value = 2469.353  # H61: 2698.6 -> 2469.353
$ws.Cells.Item(61, 9).Value = 1267.9  # I61: 1331 -> 1267.9
$ws.Cells.Item(61, 10).Value = 4185.7144  # J61: 4750 -> 4185.7144
$ws.Cells.Item(61, 11).Value = 1267.9  # K61: 1331 -> 1267.9
$ws.Cells.Item(61, 12).Value = 4185.7144  # L61: 4750 -> 4185.7144
$ws.Cells.Item(61, 13).Value = -1055.9  # M61: -1119 -> -1055.9
$ws.Cells.Item(61, 14).Value = -4609.7144  # N61: -5174 -> -4609.7144
$ws.Cells.Item(74, 8).Value = 1660.6  # H74: 1551.6666 -> 1660.6
$ws.Cells.Item(74, 9).Value = 1388.25  # I74: 1551.6666 -> 1388.25
$ws.Cells.Item(74, 10).Value = 2750  # J74: 0 -> 2750
$ws.Cells.Item(74, 11).Value = 1388.25  # K74: 1551.6666 -> 1388.25
$ws.Cells.Item(74, 12).Value = 2750  # L74: 0 -> 2750
$ws.Cells.Item(74, 13).Value = -514.25  # M74: -677.6666 -> -514.25
$ws.Cells.Item(74, 14).Value = -4498  # N74: None -> -4498
$ws.Cells.Item(77, 8).Value = 1660.6  # H77: 1551.6666 -> 1660.6
$ws.Cells.Item(77, 9).Value = 1388.25  # I77: 1551.6666 -> 1388.25
$ws.Cells.Item(77, 10).Value = 2750  # J77: 0 -> 2750
$ws.Cells.Item(77, 11).Value = 6941.25  # K77: 7758.333000000001 -> 6941.25
$ws.Cells.Item(77, 12).Value = 13750  # L77: 0 -> 13750
$ws.Cells.Item(77, 13).Value = -2573.25  # M77: -3390.333000000001 -> -2573.25
$ws.Cells.Item(77, 14).Value = -22486  # N77: None -> -22486
$ws.Cells.Item(114, 8).Value = 74399.5  # H114: 74400 -> 74399.5
$ws.Cells.Item(114, 10).Value = 74399.5  # J114: 74400 -> 74399.5
$ws.Cells.Item(114, 12).Value = 74399.5  # L114: 74400 -> 74399.5
$ws.Cells.Item(114, 14).Value = -83077.5  # N114: -83078 -> -83077.5
$ws.Cells.Item(136, 8).Value = 2469.353  # H136: 2698.6 -> 2469.353
$ws.Cells.Item(136, 9).Value = 1267.9  # I136: 1331 -> 1267.9
$ws.Cells.Item(136, 10).Value = 4185.7144  # J136: 4750 -> 4185.7144
$ws.Cells.Item(136, 11).Value = 3803.7  # K136: 3993 -> 3803.7
$ws.Cells.Item(136, 12).Value = 12557.1432  # L136: 14250 -> 12557.1432
$ws.Cells.Item(136, 13).Value = -1253.7  # M136: -1443 -> -1253.7
$ws.Cells.Item(136, 14).Value = -17657.1432  # N136: -19350 -> -17657.1432

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 2879.4  # H80: 1874.8334 -> 2879.4
$ws.Cells.Item(80, 9).Value = 132.66667  # I80: 74.5 -> 132.66667
$ws.Cells.Item(80, 10).Value = 6999.5  # J80: 2775 -> 6999.5
$ws.Cells.Item(80, 11).Value = 132.66667  # K80: 74.5 -> 132.66667
$ws.Cells.Item(80, 12).Value = 6999.5  # L80: 2775 -> 6999.5
$ws.Cells.Item(80, 13).Value = 865.3333299999999  # M80: 923.5 -> 865.3333299999999
$ws.Cells.Item(80, 14).Value = -8995.5  # N80: -4771 -> -8995.5
$ws.Cells.Item(83, 8).Value = 2879.4  # H83: 1874.8334 -> 2879.4
$ws.Cells.Item(83, 9).Value = 132.66667  # I83: 74.5 -> 132.66667
$ws.Cells.Item(83, 10).Value = 6999.5  # J83: 2775 -> 6999.5
$ws.Cells.Item(83, 11).Value = 663.3333500000001  # K83: 372.5 -> 663.3333500000001
$ws.Cells.Item(83, 12).Value = 34997.5  # L83: 13875 -> 34997.5
$ws.Cells.Item(83, 13).Value = 4328.66665  # M83: 4619.5 -> 4328.66665
$ws.Cells.Item(83, 14).Value = -44981.5  # N83: -23859 -> -44981.5
$ws.Cells.Item(134, 8).Value = 2764.8462  # H134: 8591 -> 2764.8462
$ws.Cells.Item(134, 9).Value = 1367.875  # I134: 1384 -> 1367.875
$ws.Cells.Item(134, 10).Value = 5000  # J134: 17857.143 -> 5000
$ws.Cells.Item(134, 11).Value = 4103.625  # K134: 4152 -> 4103.625
$ws.Cells.Item(134, 12).Value = 15000  # L134: 53571.429 -> 15000
$ws.Cells.Item(134, 13).Value = -1568.625  # M134: -1617 -> -1568.625
$ws.Cells.Item(134, 14).Value = -20070  # N134: -58641.429 -> -20070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 0  # H22: 400 -> 0
$ws.Cells.Item(22, 9).Value = 0  # I22: 400 -> 0
$ws.Cells.Item(22, 11).Value = 0  # K22: 400 -> 0
$ws.Cells.Item(22, 13).ClearContents()  # M22: remove (was -50)
$ws.Cells.Item(31, 8).Value = 2231.5334  # H31: 2033.7858 -> 2231.5334
$ws.Cells.Item(34, 8).Value = 2231.5334  # H34: 2033.7858 -> 2231.5334
$ws.Cells.Item(51, 8).Value = 44476  # H51: 35308.4 -> 44476
$ws.Cells.Item(51, 9).Value = 0  # I51: 30090 -> 0
$ws.Cells.Item(51, 10).Value = 44476  # J51: 38787.332 -> 44476
$ws.Cells.Item(51, 11).Value = 0  # K51: 30090 -> 0
$ws.Cells.Item(51, 12).Value = 44476  # L51: 38787.332 -> 44476
$ws.Cells.Item(51, 13).ClearContents()  # M51: remove (was -29354)
$ws.Cells.Item(51, 14).Value = -45948  # N51: -40259.332 -> -45948
$ws.Cells.Item(61, 8).Value = 44476  # H61: 35308.4 -> 44476
$ws.Cells.Item(61, 9).Value = 0  # I61: 30090 -> 0
$ws.Cells.Item(61, 10).Value = 44476  # J61: 38787.332 -> 44476
$ws.Cells.Item(61, 11).Value = 0  # K61: 30090 -> 0
$ws.Cells.Item(61, 12).Value = 44476  # L61: 38787.332 -> 44476
$ws.Cells.Item(61, 13).ClearContents()  # M61: remove (was -29742)
$ws.Cells.Item(61, 14).Value = -45172  # N61: -39483.332 -> -45172
$ws.Cells.Item(69, 8).Value = 20000  # H69: 10000 -> 20000
$ws.Cells.Item(69, 9).Value = 20000  # I69: 10000 -> 20000
$ws.Cells.Item(69, 11).Value = 20000  # K69: 10000 -> 20000
$ws.Cells.Item(69, 13).Value = -19251  # M69: -9251 -> -19251
$ws.Cells.Item(72, 8).Value = 20000  # H72: 10000 -> 20000
$ws.Cells.Item(72, 9).Value = 20000  # I72: 10000 -> 20000
$ws.Cells.Item(72, 11).Value = 60000  # K72: 30000 -> 60000
$ws.Cells.Item(72, 13).Value = -56256  # M72: -26256 -> -56256
$ws.Cells.Item(96, 8).Value = 19823.666  # H96: 20024 -> 19823.666
$ws.Cells.Item(96, 10).Value = 19823.666  # J96: 20024 -> 19823.666
$ws.Cells.Item(96, 12).Value = 19823.666  # L96: 20024 -> 19823.666
$ws.Cells.Item(96, 14).Value = -25315.666  # N96: -25516 -> -25315.666
$ws.Cells.Item(99, 8).Value = 4737.3335  # H99: 5996.25 -> 4737.3335
$ws.Cells.Item(99, 9).Value = 4069.3333  # I99: 4993 -> 4069.3333
$ws.Cells.Item(99, 10).Value = 5405.3335  # J99: 6999.5 -> 5405.3335
$ws.Cells.Item(99, 11).Value = 4069.3333  # K99: 4993 -> 4069.3333
$ws.Cells.Item(99, 12).Value = 5405.3335  # L99: 6999.5 -> 5405.3335
$ws.Cells.Item(99, 13).Value = -2571.3333  # M99: -3495 -> -2571.3333
$ws.Cells.Item(99, 14).Value = -8401.333500000001  # N99: -9995.5 -> -8401.333500000001
$ws.Cells.Item(126, 8).Value = 4737.3335  # H126: 5996.25 -> 4737.3335
$ws.Cells.Item(126, 9).Value = 4069.3333  # I126: 4993 -> 4069.3333
$ws.Cells.Item(126, 10).Value = 5405.3335  # J126: 6999.5 -> 5405.3335
$ws.Cells.Item(126, 11).Value = 12207.9999  # K126: 14979 -> 12207.9999
$ws.Cells.Item(126, 12).Value = 16216.0005  # L126: 20998.5 -> 16216.0005
$ws.Cells.Item(126, 13).Value = -9737.999899999999  # M126: -12509 -> -9737.999899999999
$ws.Cells.Item(126, 14).Value = -21156.0005  # N126: -25938.5 -> -21156.0005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 7912289.5  # H4: 8901288 -> 7912289.5
$ws.Cells.Item(4, 9).Value = 800060.2  # I4: 1000000.25 -> 800060.2
$ws.Cells.Item(4, 11).Value = 2400180.6  # K4: 3000000.75 -> 2400180.6
$ws.Cells.Item(4, 13).Value = -2400068.6  # M4: -2999888.75 -> -2400068.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(128, 8).Value = 30777  # H128: 30778 -> 30777
$ws.Cells.Item(128, 10).Value = 30777  # J128: 30778 -> 30777
$ws.Cells.Item(128, 12).Value = 30777  # L128: 30778 -> 30777
$ws.Cells.Item(128, 14).Value = -40737  # N128: -40738 -> -40737
$ws.Cells.Item(132, 8).Value = 2934.889  # H132: 2914.4285 -> 2934.889
$ws.Cells.Item(132, 9).Value = 1596.6666  # I132: 1397.25 -> 1596.6666
$ws.Cells.Item(132, 10).Value = 3604  # J132: 4937.3335 -> 3604
$ws.Cells.Item(132, 11).Value = 4789.9998  # K132: 4191.75 -> 4789.9998
$ws.Cells.Item(132, 12).Value = 10812  # L132: 14812.0005 -> 10812
$ws.Cells.Item(132, 13).Value = -2259.9998  # M132: -1661.75 -> -2259.9998
$ws.Cells.Item(132, 14).Value = -15872  # N132: -19872.0005 -> -15872

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3097.889  # H132: 3668.7144 -> 3097.889
$ws.Cells.Item(132, 9).Value = 1695.3334  # I132: 1943 -> 1695.3334
$ws.Cells.Item(132, 10).Value = 3799.1667  # J132: 4359 -> 3799.1667
$ws.Cells.Item(132, 11).Value = 5086.0002  # K132: 5829 -> 5086.0002
$ws.Cells.Item(132, 12).Value = 11397.5001  # L132: 13077 -> 11397.5001
$ws.Cells.Item(132, 13).Value = -2556.0002  # M132: -3299 -> -2556.0002
$ws.Cells.Item(132, 14).Value = -16457.5001  # N132: -18137 -> -16457.5001
$ws.Cells.Item(136, 8).Value = 3000  # H136: 2999.75 -> 3000
$ws.Cells.Item(136, 10).Value = 3000  # J136: 2999.5 -> 3000
$ws.Cells.Item(136, 12).Value = 9000  # L136: 8998.5 -> 9000
$ws.Cells.Item(136, 14).Value = -14100  # N136: -14098.5 -> -14100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 887.7273  # H113: 951.55554 -> 887.7273
$ws.Cells.Item(113, 9).Value = 837  # I113: 957.5 -> 837
$ws.Cells.Item(113, 10).Value = 1023  # J113: 939.6667 -> 1023
$ws.Cells.Item(113, 11).Value = 2511  # K113: 2872.5 -> 2511
$ws.Cells.Item(113, 12).Value = 3069  # L113: 2819.0001 -> 3069
$ws.Cells.Item(113, 13).Value = -341  # M113: -702.5 -> -341
$ws.Cells.Item(113, 14).Value = -7409  # N113: -7159.0001 -> -7409
$ws.Cells.Item(132, 8).Value = 3746  # H132: 3435.8948 -> 3746
$ws.Cells.Item(132, 9).Value = 2357.1428  # I132: 2011.1111 -> 2357.1428
$ws.Cells.Item(132, 11).Value = 7071.428400000001  # K132: 6033.3333 -> 7071.428400000001
$ws.Cells.Item(132, 13).Value = -4541.428400000001  # M132: -3503.3333 -> -4541.428400000001
$ws.Cells.Item(136, 8).Value = 1665.25  # H136: 1726.1818 -> 1665.25
$ws.Cells.Item(136, 9).Value = 1758.3  # I136: 1758.8 -> 1758.3
$ws.Cells.Item(136, 10).Value = 1200  # J136: 1400 -> 1200
$ws.Cells.Item(136, 11).Value = 5274.9  # K136: 5276.4 -> 5274.9
$ws.Cells.Item(136, 12).Value = 3600  # L136: 4200 -> 3600
$ws.Cells.Item(136, 13).Value = -2724.9  # M136: -2726.4 -> -2724.9
$ws.Cells.Item(136, 14).Value = -8700  # N136: -9300 -> -8700
